# Generate Report for Handback
#
# For the "zh-cn" and "de-de" language sheets, the handback that was
# previously recorded turned out to be out-of-scope, so the report is
# regenerated: the "Latest Target File" / "Latest Handback File" links are
# cleared, the "Latest Handback DateTime" is reset to the "never happened"
# sentinel, and an explanatory message is written into "Error Detail".
# A few report columns are also resized to better fit the (now shorter /
# longer) content.

$wb = $excel.ActiveWorkbook

$newErrorDetail = "The file with file hash 36ec3918de779b4b0f004309429f72f071337272 is not out of handoff scope."
$resetDateTime = "0001-01-01 00:00:00"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # "Latest Target File" (I2) is no longer available -- drop the value
    # and the hyperlink styling that pointed at it.
    $ws.Range("I2").Value = ""
    $ws.Range("I2").Style = "Normal"

    # "Latest Handback File" (J2) likewise no longer points at a real file.
    $ws.Range("J2").Value = ""

    # Remove the now-stale hyperlink that used to decorate I2.
    $staleHyperlinks = @()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$I$2') {
            $staleHyperlinks += $hl
        }
    }
    foreach ($hl in $staleHyperlinks) {
        $hl.Delete()
    }

    # "Latest Handback DateTime" (K2) resets to the empty/default sentinel.
    $ws.Range("K2").Value = $resetDateTime

    # "Error Detail" (P2) explains why the handback was rejected.
    $ws.Range("P2").Value = $newErrorDetail

    # Resize columns I, J and P to better match their new content.
    $ws.Columns.Item(9).ColumnWidth = 18.6506061553955 - (5 / 6)
    $ws.Columns.Item(10).ColumnWidth = 21.7054767608643 - (5 / 6)
    $ws.Columns.Item(16).ColumnWidth = 40 - (5 / 6)
}
